$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2020" column (N) mirroring the existing yearly columns (D..M).
# Row 4 holds the year headers; rows 5-17 hold the data values.

# Row 4: year header
$ws.Range("N4").Value = 2020

# Data rows 5-17 (row 15 is a blank section-header row with no values)
$ws.Range("N5").Value = 11.4
$ws.Range("N6").Value = 14.7
$ws.Range("N7").Value = 9
$ws.Range("N8").Value = 10.8
$ws.Range("N9").Value = 4.7
$ws.Range("N10").Value = 5.1
$ws.Range("N11").Value = 3.4
$ws.Range("N12").Value = 19.7
$ws.Range("N13").Value = 18.8
$ws.Range("N14").Value = 6.8
$ws.Range("N16").Value = 12.5
$ws.Range("N17").Value = 10.7

# Copy the formatting from column M (the previous last data column) onto
# column N for every row that now has data (including the blank row 15,
# which keeps its style only).
$ws.Range("M4:M17").Copy()
$ws.Range("N4:N17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the view: scroll to column E and move the active selection.
$ws.Application.ActiveWindow.ScrollColumn = 5
$sheetView = $ws
$ws.Range("A2").Select()
$ws.Application.Goto($ws.Range("E2"))
$ws.Range("S18").Select()
